$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously empty "D" values for rows 6, 7 and 20 (each row
# already has a "5" in column C; column D gets the matching "5").
$ws.Range("D6").Value = 5
$ws.Range("D7").Value = 5
$ws.Range("D20").Value = 5

# Move the active selection in the frozen (bottomRight) pane from D30 up to D5.
$ws.Range("D5").Select()
